# Updated symbol list (GitHub Actions crypto price refresh)
#
# - Refreshes the "Price" (column D) for most rank rows with the latest
#   quote. Values are stored as literal text (not numbers), matching the
#   sheet's existing convention, so a leading apostrophe is used to force
#   text entry for the numeric-looking strings.
# - A new coin (ACDXExchange) was inserted at rank 46, which pushes every
#   row below it down by one (Kangarootoken -> 47, CoinbaseStockToken -> 48,
#   BOLO -> 49, CryptobidCoin -> 50, SpecialPowerGold -> 51), and the former
#   last row (DigiFinexToken, rank 50) drops off the bottom of the sheet.
#   Column A (the rank index) does not change for these shifted rows, so
#   only B/C/D/E need to be rewritten there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Straight price (column D) refreshes -----------------------------
$ws.Range("D2").Value  = "'276.79"
$ws.Range("D3").Value  = "'21.17"
$ws.Range("D4").Value  = "'6.276"
$ws.Range("D5").Value  = "'0.06207"
$ws.Range("D6").Value  = "'3.554"
$ws.Range("D7").Value  = "'1.536"
$ws.Range("D8").Value  = "'6.582"
$ws.Range("D9").Value  = "'0.8281"
$ws.Range("D10").Value = "'0.1666"
$ws.Range("D11").Value = "'0.08302"
$ws.Range("D12").Value = "'0.03502"
$ws.Range("D13").Value = "'0.03198"
$ws.Range("D14").Value = "'0.09161"
$ws.Range("D15").Value = "'3.764"
$ws.Range("D16").Value = "'0.001631"
$ws.Range("D17").Value = "'0.04692"
$ws.Range("D18").Value = "'0.006440"
$ws.Range("D19").Value = "'0.006216"
$ws.Range("D20").Value = "'0.001067"
$ws.Range("D22").Value = "'3.720"
$ws.Range("D23").Value = "'2.316"
$ws.Range("D24").Value = "'0.01396"
$ws.Range("D40").Value = "'0.04751"
$ws.Range("D41").Value = "'0.005196"
$ws.Range("D42").Value = "'0.007068"
$ws.Range("D43").Value = "'0.1121"
$ws.Range("D44").Value = "'0.01135"
$ws.Range("D45").Value = "'0.00006320"

# --- New row inserted at rank 46 (ACDXExchange) -----------------------
$ws.Range("B46").Value = "ACDXExchange"
$ws.Range("C46").Value = "https://coinranking.com/coin/-y35lbZ7U+acdxexchange-acxt"
$ws.Range("D46").Value = "'0.0009896"
$ws.Range("E46").Value = "45ACDXExchangeACXTBestin24h"

# --- Everything below ripples down one row -----------------------------
$ws.Range("B47").Value = "Kangarootoken"
$ws.Range("C47").Value = "https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "46KangarootokenGAR"

$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.7227"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.001401"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

$ws.Range("B50").Value = "CryptobidCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D50").Value = "'0.00001899"
$ws.Range("E50").Value = "49CryptobidCoinCBC"

$ws.Range("B51").Value = "SpecialPowerGold"
$ws.Range("C51").Value = "https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
$ws.Range("D51").Value = "'0.01239"
$ws.Range("E51").Value = "50SpecialPowerGoldSPG"
